$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.123.84"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "3.501.88"
$ws.Range("E3").Value = "  +5.42%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.44"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "653.04"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.44"
$ws.Range("E7").Value = "  +6.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.417"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  +3.64%  "

$ws.Range("D11").Value = "3.503.13"
$ws.Range("E11").Value = "  +5.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.38"
$ws.Range("E12").Value = "  +11.96%  "

$ws.Range("E13").Value = "  -1.06%  "

$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.21"
$ws.Range("E14").Value = "  +4.13%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "97.008.68"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "4.155.56"
$ws.Range("E16").Value = "  +5.48%  "

$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.75"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("D19").Value = "3.516.82"
$ws.Range("E19").Value = "  +6.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.54"
$ws.Range("E20").Value = "  +11.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.04"
$ws.Range("E21").Value = "  +15.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.505"
$ws.Range("E22").Value = "  +3.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "522.73"
$ws.Range("E23").Value = "  +6.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.33"
$ws.Range("E24").Value = "  +2.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000199"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.80"
$ws.Range("E26").Value = "  +9.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "93.22"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.70"
$ws.Range("E28").Value = "  +6.58%  "

$ws.Range("D29").Value = "3.684.82"
$ws.Range("E29").Value = "  +5.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.18"
$ws.Range("E30").Value = "  +13.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("E31").Value = "  +14.30%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.141"
$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.188"
$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.41"
$ws.Range("E35").Value = "  +12.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.590"
$ws.Range("E36").Value = "  +9.10%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.96"
$ws.Range("E38").Value = "  +6.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.154"
$ws.Range("E40").Value = "  +3.80%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "515.40"
$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.921"
$ws.Range("E43").Value = "  +12.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.30"
$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0427"
$ws.Range("E45").Value = "  +5.71%  "

$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  +7.05%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.44"
$ws.Range("E47").Value = "  +10.95%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.64"
$ws.Range("E48").Value = "  +4.59%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.23"
$ws.Range("E49").Value = "  +14.38%  "

$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.60"
$ws.Range("E50").Value = "  -2.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.53"
$ws.Range("E51").Value = "  +2.07%  "
